$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate the content of rows 8, 9 and 10 ---
# Before: row8=json/pattern, row9=int/smallint-max, row10=float/number
# After : row8=(old row9), row9=(old row10), row10=(old row8)
$b8 = $ws.Range("B8").Value2
$c8 = $ws.Range("C8").Value2
$b9 = $ws.Range("B9").Value2
$c9 = $ws.Range("C9").Value2
$b10 = $ws.Range("B10").Value2
$c10 = $ws.Range("C10").Value2

$ws.Range("B8").Value2 = $b9
$ws.Range("C8").Value2 = $c9

$ws.Range("B9").Value2 = $b10
$ws.Range("C9").Value2 = $c10

$ws.Range("B10").Value2 = $b8
$ws.Range("C10").Value2 = $c8

# --- Row heights follow the content that moved into each row ---
$ws.Rows.Item(8).RowHeight = 38.25
$ws.Rows.Item(9).EntireRow.AutoFit()
$ws.Rows.Item(10).RowHeight = 63.75

# --- Update the frozen-pane view: top-left visible cell moves from C4 to C3,
#     and the active selection in the bottom-right pane moves from C9 to C3 ---
$win = $excel.ActiveWindow
$ws.Range("C3").Select()
$win.ScrollRow = 3
$win.ScrollColumn = 3
